$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B (shifts B..J -> C..K)
$ws.Columns.Item(2).Insert()

# 2. Bold header label in A3 ("lineare SVM")
$ws.Range("A3").Font.Bold = $true

# 3. Add "1-Accuracy" helper formulas in column K for the first table
$ws.Range("K4").Formula = "=1-H4"
$ws.Range("K5").Formula = "=1-H5"
$ws.Range("K6").Formula = "=1-H6"

# 4. Second table header (row 8): Execution / Accuracy / TPR / FPR / Classification Error
#    Re-use the same header style as row 3 (G3:K3) by copying its formatting.
$ws.Range("G3:K3").Copy()
$ws.Range("G8:K8").PasteSpecial(-4122)
$ws.Range("G8").Value = "Execution"
$ws.Range("H8").Value = "Accuracy"
$ws.Range("I8").Value = "TPR"
$ws.Range("J8").Value = "FPR"
$ws.Range("K8").Value = "Classification Error"
$ws.Range("G8:K8").RowHeight = 30

# 5. Execution numbers 1..5 down column G (rows 9-13), with wrap text applied down to row 51
$ws.Range("G9:G51").WrapText = $true
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 3
$ws.Range("G12").Value = 4
$ws.Range("G13").Value = 5

# 6. Explicit zero accuracy/TPR/FPR values for execution 2
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0

# 7. Classification Error column for executions 1-5 (shared formula =1-H)
$ws.Range("K9").Formula = "=1-H9"
$ws.Range("K10").Formula = "=1-H10"
$ws.Range("K11").Formula = "=1-H11"
$ws.Range("K12").Formula = "=1-H12"
$ws.Range("K13").Formula = "=1-H13"

# 8. Mean / std summary block to the right (M11:N13)
$ws.Range("M11:N11").Merge()
$ws.Range("M11").Value = "accuracy"
$ws.Range("M11").HorizontalAlignment = -4108
$ws.Range("M12").Value = "mean"
$ws.Range("N12").Value = "std"
$ws.Range("M13").Formula = "=AVERAGE(H12:H16)"
$ws.Range("M13").NumberFormat = "0.0000000"
$ws.Range("N13").Formula = "=_xlfn.STDEV.P(H12:H16)"

# 9. Second classifier label, bold, at A16
$ws.Range("A16").Value = "rbf SVM"
$ws.Range("A16").Font.Bold = $true

# 10. Recalculate so cached values are correct
$ws.Calculate()
